# Weekly update: insert a new latest-week record at the top of the data
# (row 2), pushing all existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (first data row).
$ws.Rows("2:2").Insert()

# Excel inherits formatting from the row above (the bold header row) when
# inserting, so reset the new row's formatting and re-apply the same
# number format used by the other data rows in column D (the date column).
$ws.Range("A2:R2").ClearFormats()
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# Populate the new row with the latest week's data.
$ws.Range("A2").Value2 = 1
$ws.Range("B2").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value2 = "Arica y Parinacota"
$ws.Range("D2").Value2 = 44847
$ws.Range("E2").Value2 = 15
$ws.Range("F2").Value2 = 100112003
$ws.Range("G2").Value2 = "Ajo"
$ws.Range("H2").Value2 = "Chino"
$ws.Range("I2").Value2 = "Primera"
$ws.Range("J2").Value2 = 400
$ws.Range("K2").Value2 = 16000
$ws.Range("L2").Value2 = 17000
$ws.Range("M2").Value2 = 16500
$ws.Range("N2").Value2 = "`$/caja 10 kilos"
$ws.Range("O2").Value2 = "China"
$ws.Range("P2").Value2 = 1650
$ws.Range("Q2").Value2 = 10
$ws.Range("R2").Value2 = "Hortaliza"
